$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "Test - Test - 10/30/2020" values to "Test - Test"
$ws.Range("D2").Value = "Test - Test"
$ws.Range("D4").Value = "Test - Test"

# Update selection to E4 (single cell) instead of A4:H4
$ws.Range("E4").Select()
